# Add a new benchmark result row for the 12-core Opteron 6180 SE.
# The row is inserted as row 10 (keeping the table's ascending-by-column-O
# sort order), pushing the existing rows 10-22 down to 11-23.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 10 (shifts rows 10..22 down to 11..23,
# and extends the table / used range accordingly).
$ws.Rows.Item(10).Insert()

# Fill in the new row with the Opteron 6180 SE benchmark data.
$ws.Range("A10").Value = "AMD"
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = "Opteron 6180 SE"
$ws.Range("D10").Value = 140
$ws.Range("E10").Value = 12
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2.5
$ws.Range("H10").Value = 2.5
$ws.Range("I10").Value = "x86-64"
$ws.Range("J10").Value = "Performance"
$ws.Range("K10").Value = 32
$ws.Range("L10").Value = 2
$ws.Range("M10").Value = "DDR3"
$ws.Range("N10").Value = 1333
$ws.Range("O10").Value = 1.1200000000000001
$ws.Range("P10").Value = 2.17
$ws.Range("Q10").Value = 4.16
$ws.Range("R10").Value = 8.32

# Match the final selection left by the author in the saved file.
$ws.Range("C10").Select()
